# Update Iceland MSME country indicator values with more precise figures,
# while preserving the original (General) cell formatting/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 13: Enterprises density (per 1000 people)
Set-TextValue "B13" "74.55"
Set-TextValue "C13" "4.64"

# Row 14: Employment (% of total)
Set-TextValue "B14" "26.31"
Set-TextValue "C14" "45.36"
Set-TextValue "D14" "71.67"

# Row 16: Enterprises (% of total)
Set-TextValue "B16" "93.91"
Set-TextValue "C16" "5.85"
Set-TextValue "D16" "99.76"
